$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.918.64'
$ws.Cells.Item(2, 5).Value = '  +0.32%  '
$ws.Cells.Item(3, 4).Value = '2.535.07'
$ws.Cells.Item(3, 5).Value = '  -0.84%  '
$ws.Cells.Item(4, 4).Value = "'" + '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).Value = "'" + '311.03'
$ws.Cells.Item(5, 5).Value = '  -0.01%  '
$ws.Cells.Item(6, 4).Value = "'" + '100.68'
$ws.Cells.Item(6, 5).Value = '  +2.32%  '
$ws.Cells.Item(7, 4).Value = "'" + '0.568'
$ws.Cells.Item(7, 5).Value = '  -0.64%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 5).Value = '  -1.25%  '
$ws.Cells.Item(10, 4).Value = "'" + '35.82'
$ws.Cells.Item(10, 5).Value = '  +0.61%  '
$ws.Cells.Item(11, 4).Value = "'" + '0.0807'
$ws.Cells.Item(11, 5).Value = '  -0.15%  '
$ws.Cells.Item(12, 5).Value = '  -1.06%  '
$ws.Cells.Item(13, 5).Value = '  +1.59%  '
$ws.Cells.Item(14, 4).Value = '2.924.37'
$ws.Cells.Item(14, 5).Value = '  -0.88%  '
$ws.Cells.Item(15, 4).Value = "'" + '15.49'
$ws.Cells.Item(15, 5).Value = '  -3.05%  '
$ws.Cells.Item(16, 4).Value = '2.544.40'
$ws.Cells.Item(16, 5).Value = '  -0.45%  '
$ws.Cells.Item(17, 4).Value = "'" + '0.818'
$ws.Cells.Item(17, 5).Value = '  -2.36%  '
$ws.Cells.Item(18, 4).Value = '42.882.72'
$ws.Cells.Item(18, 5).Value = '  +0.21%  '
$ws.Cells.Item(19, 4).Value = "'" + '6.70'
$ws.Cells.Item(19, 5).Value = '  -0.55%  '
$ws.Cells.Item(20, 4).Value = "'" + '12.40'
$ws.Cells.Item(20, 5).Value = '  +0.34%  '
$ws.Cells.Item(21, 5).Value = '  -0.47%  '
$ws.Cells.Item(22, 5).Value = '  +0.35%  '
$ws.Cells.Item(23, 4).Value = "'" + '244.05'
$ws.Cells.Item(23, 5).Value = '  -1.38%  '
$ws.Cells.Item(24, 4).Value = "'" + '2.88'
$ws.Cells.Item(24, 5).Value = '  -1.26%  '
$ws.Cells.Item(25, 4).Value = "'" + '2.05'
$ws.Cells.Item(25, 5).Value = '  +0.32%  '
$ws.Cells.Item(26, 5).Value = '  +0.06%  '
$ws.Cells.Item(27, 4).Value = "'" + '25.64'
$ws.Cells.Item(27, 5).Value = '  -4.25%  '
$ws.Cells.Item(28, 5).Value = '  -2.64%  '
$ws.Cells.Item(29, 4).Value = "'" + '10.24'
$ws.Cells.Item(29, 5).Value = '  +0.67%  '
$ws.Cells.Item(30, 4).Value = "'" + '38.91'
$ws.Cells.Item(30, 5).Value = '  -2.70%  '
$ws.Cells.Item(31, 4).Value = "'" + '160.24'
$ws.Cells.Item(31, 5).Value = '  +1.21%  '
$ws.Cells.Item(32, 4).Value = "'" + '5.83'
$ws.Cells.Item(32, 5).Value = '  +1.48%  '
$ws.Cells.Item(33, 5).Value = '  +7.62%  '
$ws.Cells.Item(34, 5).Value = '  -0.23%  '
$ws.Cells.Item(35, 5).Value = '  +1.23%  '
$ws.Cells.Item(36, 4).Value = "'" + '18.37'
$ws.Cells.Item(36, 5).Value = '  -1.56%  '
$ws.Cells.Item(37, 5).Value = '  -3.95%  '
$ws.Cells.Item(38, 5).Value = '  -4.80%  '
$ws.Cells.Item(39, 5).Value = '  +0.06%  '
$ws.Cells.Item(40, 5).Value = '  +0.18%  '
$ws.Cells.Item(41, 4).Value = "'" + '4.19'
$ws.Cells.Item(41, 5).Value = '  +3.13%  '
$ws.Cells.Item(42, 4).Value = "'" + '21.99'
$ws.Cells.Item(42, 5).Value = '  -2.89%  '
$ws.Cells.Item(43, 4).Value = "'" + '3.34'
$ws.Cells.Item(43, 5).Value = '  +4.53%  '
$ws.Cells.Item(44, 5).Value = '  +0.22%  '
$ws.Cells.Item(45, 5).Value = '  -0.39%  '
$ws.Cells.Item(46, 4).Value = '2.003.91'
$ws.Cells.Item(46, 5).Value = '  +0.74%  '
$ws.Cells.Item(47, 4).Value = "'" + '9.30'
$ws.Cells.Item(47, 5).Value = '  +3.72%  '
$ws.Cells.Item(48, 4).Value = '2.775.33'
$ws.Cells.Item(49, 4).Value = "'" + '0.193'
$ws.Cells.Item(49, 5).Value = '  -0.21%  '
$ws.Cells.Item(50, 4).Value = "'" + '79.90'
$ws.Cells.Item(50, 5).Value = '  -1.70%  '
$ws.Cells.Item(51, 4).Value = "'" + '72.59'
$ws.Cells.Item(51, 5).Value = '  -1.13%  '
